$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the Tue Jul 25 11:33:16 UTC 2023 GitHub Actions cryptos data
# refresh: updated Price (D) / Volume 1h (E) figures for each listed
# coin, plus the FraxShare <-> TheSandbox row swap (rows 43-44).
#
# Column D holds plain text (not numbers) in the source data, even
# when the text is numeric-looking (e.g. "237.82"). A leading
# apostrophe forces Excel to store it as text instead of coercing
# it to a Double, matching the original inline-string cell type.

$ws.Range("D2").Value = '29.166.63'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.853.68'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''237.82'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = '''0.6875'
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '''0.07782'
$ws.Range("E8").Value = '  +3.00%  '
$ws.Range("D9").Value = '''0.3037'
$ws.Range("E9").Value = '  -1.45%  '
$ws.Range("D10").Value = '''23.16'
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("D11").Value = '''0.08086'
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '1.840.28'
$ws.Range("E12").Value = '  -1.39%  '
$ws.Range("D13").Value = '''0.7200'
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("D14").Value = '''5.189'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '''89.20'
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").Value = '29.169.48'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '''5.732'
$ws.Range("E17").Value = '  -3.00%  '
$ws.Range("D18").Value = '''0.000007795'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").Value = '''13.27'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '''234.04'
$ws.Range("E20").Value = '  -3.67%  '
$ws.Range("D21").Value = '''1.0000'
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").Value = '2.109.79'
$ws.Range("E22").Value = '  -0.79%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '''7.475'
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").Value = '''161.57'
$ws.Range("E25").Value = '  -0.30%  '
$ws.Range("D26").Value = '''8.961'
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("D27").Value = '''0.1429'
$ws.Range("E27").Value = '  -2.94%  '
$ws.Range("D28").Value = '''18.03'
$ws.Range("E28").Value = '  -0.38%  '
$ws.Range("D29").Value = '''1.955'
$ws.Range("E29").Value = '  +0.80%  '
$ws.Range("D30").Value = '''1.404'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '''4.496'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D32").Value = '''1.480'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = '''4.007'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("D34").Value = '''0.05205'
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("D35").Value = '''1.177'
$ws.Range("E35").Value = '  -1.88%  '
$ws.Range("D36").Value = '''0.7019'
$ws.Range("E36").Value = '  -2.45%  '
$ws.Range("D37").Value = '''1.002'
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").Value = '''2.674'
$ws.Range("E38").Value = '  +0.39%  '
$ws.Range("D39").Value = '''0.01844'
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("D40").Value = '''2.692'
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("D41").Value = '''0.9334'
$ws.Range("E41").Value = '  +5.53%  '
$ws.Range("D42").Value = '1.101.62'
$ws.Range("E42").Value = '  +5.30%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '''0.4280'
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '''5.911'
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = '''70.30'
$ws.Range("E45").Value = '  +0.24%  '
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").Value = '''102.47'
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("D48").Value = '''1.791'
$ws.Range("E48").Value = '  +2.26%  '
$ws.Range("D49").Value = '2.005.48'
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("D50").Value = '''9.153'
$ws.Range("E50").Value = '  -1.30%  '
$ws.Range("D51").Value = '''6.991'
$ws.Range("E51").Value = '  -3.88%  '
